$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. thousands separated by
# dots, like "35.226.62"), stored as text in the original workbook. Force the
# cells we touch in column D to keep a Text number format so Excel does not
# reinterpret the assigned strings as numbers (which would corrupt values like
# "35.226.62" or introduce float rounding such as "245.28999999999999").

$dCells = @("D2","D3","D5","D8","D10","D11","D13","D14","D15","D16","D17","D18","D19","D21","D22","D23","D26","D27","D28","D29","D32","D33","D34","D37","D39","D41","D42","D44","D45","D49","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "35.237.74"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "1.886.83"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "245.29"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "42.70"
$ws.Range("E8").Value = "  +3.67%  "

$ws.Range("E9").Value = "  +3.44%  "

$ws.Range("D10").Value = "54.81"
$ws.Range("E10").Value = "  +7.29%  "

$ws.Range("D11").Value = "0.0740"
$ws.Range("E11").Value = "  +2.03%  "

$ws.Range("E12").Value = "  +1.85%  "

$ws.Range("D13").Value = "13.73"
$ws.Range("E13").Value = "  +8.04%  "

$ws.Range("D14").Value = "0.777"
$ws.Range("E14").Value = "  +10.52%  "

$ws.Range("D15").Value = "2.158.89"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "4.95"
$ws.Range("E16").Value = "  +3.09%  "

$ws.Range("D17").Value = "1.879.69"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "35.259.88"
$ws.Range("E18").Value = "  +1.61%  "

$ws.Range("D19").Value = "73.12"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("E20").Value = "  +1.87%  "

$ws.Range("D21").Value = "243.53"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").Value = "12.76"
$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  +6.42%  "

$ws.Range("E24").Value = "  +7.77%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D28").Value = "8.48"
$ws.Range("E28").Value = "  +2.04%  "

$ws.Range("D29").Value = "18.22"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("D32").Value = "0.0592"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  -13.94%  "

$ws.Range("D37").Value = "0.848"
$ws.Range("E37").Value = "  +3.91%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("D39").Value = "0.0715"
$ws.Range("E39").Value = "  +7.58%  "

$ws.Range("E40").Value = "  +5.30%  "

$ws.Range("D41").Value = "97.85"
$ws.Range("E41").Value = "  +1.51%  "

$ws.Range("D42").Value = "17.07"
$ws.Range("E42").Value = "  +1.69%  "

$ws.Range("E43").Value = "  +0.85%  "

$ws.Range("D44").Value = "1.323.91"
$ws.Range("E44").Value = "  +4.00%  "

$ws.Range("D45").Value = "13.18"
$ws.Range("E45").Value = "  +11.47%  "

$ws.Range("E46").Value = "  +2.41%  "

$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").Value = "2.74"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").Value = "2.059.18"
$ws.Range("E51").Value = "  +0.80%  "

# Row reorderings: Monero/PancakeSwap (26<->27) and WEMIXToken/InternetComputer(DFINITY) (33<->34)
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.14"
$ws.Range("E26").Value = "  -1.31%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "167.21"
$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.16"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +16.89%  "
